$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "'1"
$ws.Range("D2").Value = "'01-01-2024"
$ws.Range("E2").Value = "'31-01-2024"

$ws.Range("A3").Value = "D21CQCN01-N"
$ws.Range("B3").Value = "2023-2024"
$ws.Range("C3").Value = "'2"
$ws.Range("D3").Value = "'01-03-2024"
$ws.Range("E3").Value = "'15-06-2024"

$ws.Range("A4").Value = "D22CQCN01-N"
$ws.Range("B4").Value = "2023-2024"
$ws.Range("C4").Value = "'2"
$ws.Range("D4").Value = "'25-05-2024"
$ws.Range("E4").Value = "'25-06-2024"

$ws.Range("A1:E4").ClearFormats()
